# BOM_miniFOC.xlsx edit: "The PCB circuit of minifoc is modified"
#  - J1,J2 connector entry (row 6) becomes a standalone "J2" entry
#  - L1 inductor (row 8) replaced by L2 (new package/LCSC part/price)
#  - A new row 19 is added for connector J5 (GH1.25 3-pin)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: L1 -> L2 inductor, new footprint / LCSC part / unit price ---
$ws.Cells.Item(8, 1).Value = "L2"
$ws.Cells.Item(8, 5).Value = "C251690"
$ws.Cells.Item(8, 3).Value = "SMD_L4.4-W4.0"
$ws.Cells.Item(8, 6).Value = 1.9060999999999999

# --- Row 19 (new): connector J5, GH1.25 3-pin ---
$ws.Cells.Item(19, 1).Value = "J5"
$ws.Cells.Item(19, 5).Value = "C2829252"

# --- Row 6: designator J1,J2 -> J2 (part/package/qty/price unchanged) ---
$ws.Cells.Item(6, 1).Value = "J2"

$ws.Cells.Item(19, 2).Value = "GH1.25"
$ws.Cells.Item(19, 3).Value = "GH1.25_3P_卧贴"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3221

# Match the formatting of the other data rows (centered text/values with the
# thin border used throughout the table) for the newly added row.
$newRow = $ws.Range("A19:F19")
$newRow.Borders.LineStyle = 1
$newRow.HorizontalAlignment = -4108

# Update selection to reflect where the author left off editing.
$ws.Range("H15").Select()
